$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AX. This shifts the existing "nom" column
# (AX -> AY) and "url_produit" column (AY -> AZ) one position to the right,
# and grows the used range from A1:AY206 to A1:AZ206.
$ws.Columns("AX:AX").Insert()

# The new AX column is a fresh price snapshot column: it carries the same
# values (or blanks) as the previous last snapshot column (AW), mirroring
# how every prior snapshot column was appended.
$ws.Range("AW2:AW206").Copy($ws.Range("AX2:AX206"))

# AX1 gets the new snapshot's timestamp header (same style as the other
# header cells in row 1).
$ws.Range("AX1").Value = "2026-01-29 21:16:42"
